$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '44.133.77'
$ws.Range('E2').Value = '  +4.67%  '
$ws.Range('D3').Value = '2.225.16'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '82.85'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +13.39%  '
$ws.Range('E7').Value = '  +3.31%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.612'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.28'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +11.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0937'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.08'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.39%  '
$ws.Range('E13').Value = '  +2.96%  '
$ws.Range('D14').Value = '2.563.94'
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.69'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.88%  '
$ws.Range('D16').Value = '2.235.95'
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.786'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.24%  '
$ws.Range('D18').Value = '44.008.57'
$ws.Range('E18').Value = '  +4.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000105'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.71'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('E21').Value = '  +4.09%  '
$ws.Range('E22').Value = '  +11.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.66'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.49%  '
$ws.Range('E24').Value = '  -2.99%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.80'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '40.76'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +12.04%  '
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.53'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0898'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +11.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.69'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.38'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.55%  '
$ws.Range('E35').Value = '  +10.11%  '
$ws.Range('E36').Value = '  +2.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0365'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.53'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.53'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +14.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.01'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +26.52%  '
$ws.Range('E41').Value = '  +4.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.14'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +9.74%  '
$ws.Range('E43').Value = '  +8.22%  '
$ws.Range('E44').Value = '  +4.30%  '
$ws.Range('E45').Value = '  +1.89%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.41'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0988'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.68%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.13'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.57%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.57'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +29.98%  '
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('E51').Value = '  +4.10%  '
